$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new test-result rows (76 and 77) to the Campus Test Results sheet,
# mirroring the existing rows' layout: A=UUID, B=Test name, C=Status,
# D=Start time, E=End time, F=Duration.

$ws.Range("A76").Value = "65840174-3784-43b7-86fe-cfe855126c5e"
$ws.Range("B76").Value = "Login with valid credentials"
$ws.Range("C76").Value = "PASSED"
$ws.Range("D76").Value = "03_29_2024_00_15_59"
$ws.Range("E76").Value = "03_29_2024_00_16_05"
$ws.Range("F76").Value = "PT5.90775S"

$ws.Range("A77").Value = "c11cb73a-64df-47c9-be4c-b157d0f70c74"
$ws.Range("B77").Value = "Create Country"
$ws.Range("C77").Value = "PASSED"
$ws.Range("D77").Value = "03_29_2024_00_16_09"
$ws.Range("E77").Value = "03_29_2024_00_16_17"
$ws.Range("F77").Value = "PT8.3585238S"
